$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 690 (everything from row 690 down shifts by 2,
# e.g. old row 690 -> new row 692, old row 731 -> new row 733).
$ws.Range("A690:A691").EntireRow.Insert()

# Fill in the two freshly inserted rows with the new data points.
# Column A holds plain text dates (e.g. "2026/01/22"), so force a text
# number format before writing them to stop them being auto-converted into
# Excel date serial numbers, then clear the formatting again afterwards so
# the cell is left with no explicit style (matching the rest of the sheet).
$newRows = @(
    @(690, "2026/01/22", "木", 7),
    @(691, "2026/01/22", "木", 10)
)

foreach ($item in $newRows) {
    $r = $item[0]
    $dateText = $item[1]
    $dow = $item[2]
    $hour = $item[3]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dateText

    $ws.Cells.Item($r, 2).Value = $dow
    $ws.Cells.Item($r, 3).Value = $hour
    $ws.Cells.Item($r, 4).Value = 201
}

$ws.Range("A690:A691").ClearFormats()
